# "remove column from alcohol data"
# The measurement sheet (Sheet1) had a duplicate/extra column (M) whose data
# was redundant with the following column (N). The fix deletes column M
# entirely, which shifts the old column N (and everything after it, though
# there was nothing else) one place to the left so it becomes the new
# column M. This is a straightforward "delete entire column" operation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete column M (the 13th column) - this removes its values and shifts
# column N (now holding the surviving data) left into column M.
$ws.Columns.Item(13).Delete() | Out-Null

# Leave the selection on the cell that now represents the boundary of the
# trimmed data range, matching where Excel would naturally leave the cursor
# after removing the column.
$ws.Range("M1").Select() | Out-Null
